$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param([string]$CellRef, [string]$NewValue)
    $r = $ws.Range($CellRef)
    $r.NumberFormat = "@"
    $r.Value2 = $NewValue
    $r.Style = "Normal"
}

Set-TextCell "D2" "58.889.47"
Set-TextCell "E2" "  -2.83%  "
Set-TextCell "D3" "2.634.78"
Set-TextCell "E3" "  -2.44%  "
Set-TextCell "E4" "  -0.06%  "
Set-TextCell "D5" "526.32"
Set-TextCell "E5" "  -0.17%  "
Set-TextCell "D6" "144.03"
Set-TextCell "E6" "  -3.71%  "
Set-TextCell "E7" "  +0.10%  "
Set-TextCell "D8" "0.570"
Set-TextCell "E8" "  -1.32%  "
Set-TextCell "D9" "6.63"
Set-TextCell "E9" "  -5.85%  "
Set-TextCell "E10" "  -1.53%  "
Set-TextCell "D11" "0.336"
Set-TextCell "E11" "  -1.55%  "
Set-TextCell "E12" "  +0.62%  "
Set-TextCell "D13" "3.099.47"
Set-TextCell "E13" "  -2.44%  "
Set-TextCell "D14" "58.861.73"
Set-TextCell "E14" "  -2.88%  "
Set-TextCell "D15" "20.97"
Set-TextCell "E15" "  -2.50%  "
Set-TextCell "D16" "0.0000137"
Set-TextCell "E16" "  -1.37%  "
Set-TextCell "D17" "2.640.67"
Set-TextCell "E17" "  -2.51%  "
Set-TextCell "B18" "Polkadot"
Set-TextCell "C18" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextCell "D18" "4.46"
Set-TextCell "E18" "  -1.05%  "
Set-TextCell "B19" "BitcoinCash"
Set-TextCell "C19" "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextCell "D19" "340.46"
Set-TextCell "E19" "  -1.70%  "
Set-TextCell "D20" "10.53"
Set-TextCell "E20" "  -0.13%  "
Set-TextCell "D21" "6.33"
Set-TextCell "E21" "  -1.09%  "
Set-TextCell "E22" "  +0.26%  "
Set-TextCell "D23" "65.44"
Set-TextCell "E23" "  +2.84%  "
Set-TextCell "E24" "  -0.37%  "
Set-TextCell "E25" "  -2.03%  "
Set-TextCell "D26" "0.998"
Set-TextCell "E26" "  +0.23%  "
Set-TextCell "D27" "7.23"
Set-TextCell "E27" "  -0.86%  "
Set-TextCell "D28" "0.0₃0797"
Set-TextCell "E28" "  -3.07%  "
Set-TextCell "D29" "6.48"
Set-TextCell "E29" "  -4.56%  "
Set-TextCell "E30" "  -0.02%  "
Set-TextCell "E31" "  +0.33%  "
Set-TextCell "D32" "18.90"
Set-TextCell "E32" "  -0.94%  "
Set-TextCell "D33" "149.79"
Set-TextCell "E33" "  -0.63%  "
Set-TextCell "D34" "4.21"
Set-TextCell "E34" "  -1.42%  "
Set-TextCell "D35" "1.20"
Set-TextCell "E35" "  -2.15%  "
Set-TextCell "D36" "0.922"
Set-TextCell "E36" "  +0.02%  "
Set-TextCell "D37" "0.872"
Set-TextCell "E37" "  -3.52%  "
Set-TextCell "D38" "36.57"
Set-TextCell "E38" "  -2.08%  "
Set-TextCell "E39" "  -3.84%  "
Set-TextCell "D40" "3.65"
Set-TextCell "E40" "  -0.25%  "
Set-TextCell "D41" "0.997"
Set-TextCell "E41" "  +0.13%  "
Set-TextCell "B42" "Stellar"
Set-TextCell "C42" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextCell "D42" "0.0975"
Set-TextCell "E42" "  -1.10%  "
Set-TextCell "B43" "Mantle"
Set-TextCell "C43" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextCell "D43" "0.603"
Set-TextCell "E43" "  -5.08%  "
Set-TextCell "D44" "270.52"
Set-TextCell "E44" "  -3.31%  "
Set-TextCell "B45" "EnergySwap"
Set-TextCell "C45" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D45" "19.35"
Set-TextCell "E45" "  -3.80%  "
Set-TextCell "D46" "0.0538"
Set-TextCell "E46" "  -1.10%  "
Set-TextCell "B47" "WhiteBITCoin"
Set-TextCell "C47" "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextCell "D47" "10.65"
Set-TextCell "E47" "  +1.00%  "
Set-TextCell "D48" "2.040.57"
Set-TextCell "E48" "  -2.74%  "
Set-TextCell "B49" "VeChain"
Set-TextCell "C49" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell "D49" "0.0230"
Set-TextCell "E49" "  -1.56%  "
Set-TextCell "B50" "InjectiveProtocol"
Set-TextCell "C50" "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextCell "D50" "19.03"
Set-TextCell "E50" "  -2.29%  "
Set-TextCell "B51" "RenderToken"
Set-TextCell "C51" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell "D51" "4.68"
Set-TextCell "E51" "  -6.13%  "
